$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.882.25"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.07%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.198.68"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.13%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "596.21"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.32%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.69"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.17%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.195.88"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.11%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.533"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.88%  "
$ws.Range("E10").Value = "  -0.63%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.07"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.64%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.513"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.41%  "
$ws.Range("E13").Value = "  +1.18%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "38.95"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.30%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.722.34"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.03%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.901.14"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.96%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.41"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.67%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.200.77"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.13%  "
$ws.Range("E19").Value = "  -0.03%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "509.18"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.69%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.31"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.94%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.741"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.77%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "15.25"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.64%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.71%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.86"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.04%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.998"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.30%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.33"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.67%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.99"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.10%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.26"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.51%  "
$ws.Range("E30").Value = "  +1.36%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.86"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +8.60%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "28.23"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.51%  "
$ws.Range("E33").Value = "  +2.55%  "
$ws.Range("E34").Value = "  +0.02%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.55"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.26%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "55.02"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.45%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0906"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.17%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "481.56"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.79%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0419"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.10%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.92"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.19%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.83"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.76%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.298"
$ws.Range("D42").Style = "Normal"
$ws.Range("E43").Value = "  +2.64%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0₃0646"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +8.27%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.934.76"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.29%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.42"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.34%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "28.27"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.21%  "
$ws.Range("E49").Value = "  +1.15%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.29"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.45%  "
$ws.Range("B51").Value = "CoreDAO"
$ws.Range("C51").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.54"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.44%  "
